# ============================================================================
# Applies the "Add files via upload" revision to conditions.xlsx:
#  - Re-randomises the color/target/corrAns/sound columns on several trial
#    rows, switches the "0.1" duration rows from text to numeric 0.1,
#    clears a handful of stray values, and appends three new trial rows
#    (43-45).
#  - Styles the header row (A1:F1) with bold text, a thin box border and
#    centered/top aligned text.
#  - Updates the window view (zoom / selection).
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Header styling (A1:F1): bold font, thin box border, centered + top
#    aligned. Build the format on a scratch cell first and copy it across
#    in a single PasteSpecial so only one new font/border/style entry is
#    created (instead of one per incremental property change).
# ----------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160

$scratch.Copy() | Out-Null
$headerRange = $ws.Range("A1:F1")
$headerRange.PasteSpecial(-4122) | Out-Null

$scratch.ClearContents()

# ----------------------------------------------------------------------
# 2. Trial-data edits (rows 2-45)
# ----------------------------------------------------------------------
$ws.Range("C2").Value = 0.1

$ws.Range("B14").Value = "green.png"
$ws.Range("C14").Value = 0.1
$ws.Range("F14").Value = "silent.wav"

$ws.Range("B15").Value = "orange.png"
$ws.Range("F15").Value = "beep.wav"

$ws.Range("B18").Value = "green.png"
$ws.Range("F18").Value = "silent.wav"

$ws.Range("B19").Value = "orange.png"
$ws.Range("D19").Value = 0

$ws.Range("C20").Value = 0.1

$ws.Range("B21").Value = "red.png"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "beep.wav"

$ws.Range("B23").Value = "green.png"
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = "silent.wav"

$ws.Range("B25").Value = "red.png"
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = "beep.wav"

$ws.Range("B26").Value = "green.png"

$ws.Range("F28").Value = "silent.wav"

$ws.Range("B30").Value = "orange.png"
$ws.Range("D30").Value = 0

$ws.Range("B32").Value = "red.png"
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 1

$ws.Range("B34").Value = "orange.png"
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = "beep.wav"

$ws.Range("B35").Value = "green.png"
$ws.Range("D35").Value = 0
$ws.Range("F35").Value = "silent.wav"

$ws.Range("B37").Value = "red.png"
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = "beep.wav"

$ws.Range("B39").Value = "green.png"
$ws.Range("F39").Value = "silent.wav"

$ws.Range("B40").Value = "green.png"
$ws.Range("D40").Value = 0
$ws.Range("F40").Value = "silent.wav"

$ws.Range("B41").Value = "orange.png"
$ws.Range("F41").Value = "beep.wav"

$ws.Range("C42").Value = 0.1

# New trial rows appended at the bottom
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "red.png"
$ws.Range("C43").Value = 3
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = "beep.wav"

$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "green.png"
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = "silent.wav"

$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "green.png"
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = "silent.wav"

# ----------------------------------------------------------------------
# 3. Clear stray cells that must end up empty
# ----------------------------------------------------------------------
$ws.Range("A14").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("A20").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("E35").ClearContents()
$ws.Range("E39").ClearContents()
$ws.Range("E40").ClearContents()
$ws.Range("A42").ClearContents()
$ws.Range("E42").ClearContents()

# ----------------------------------------------------------------------
# 4. Recalculate (so H5 = SUM(C:C) reflects the new 0.1 numeric values)
# ----------------------------------------------------------------------
$excel.Calculate() | Out-Null

# ----------------------------------------------------------------------
# 5. Window / view state
# ----------------------------------------------------------------------
$ws.Range("C42").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 119
$win.ScrollRow = 36
$win.ScrollColumn = 1
